{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst VT = \"\\u000b\"; // <w:br/> maps to a vertical-tab char in Office.js text\n\n// Paragraph 0: title \"LTI Solution\" -> \"LTI - Java Solution\"\nparagraphs.items[0].insertText(\"LTI - Java Solution\", \"Replace\");\n\n// Paragraph 1: question 1 - change text and add a new \"b)w\" answer line\nparagraphs.items[1].insertText(\n  \"1)adf\" + VT + \"a)d\" + VT + \"b)w\" + VT + VT,\n  \"Replace\"\n);\n\n// Paragraph 2: question 2 - change text and remove the \"b)3\" answer line\nparagraphs.items[2].insertText(\n  \"2)ffff\" + VT + \"a)1\" + VT + VT,\n  \"Replace\"\n);\n\n// Paragraph 3: question 3 - change text and remove the \"b)o\" answer line\nparagraphs.items[3].insertText(\n  \"3)dfa\" + VT + \"a)d\" + VT + VT,\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$vt = [char]11\n\n# Paragraph 1: title \"LTI Solution\" -> \"LTI - Java Solution\"\n$d.Paragraphs(1).Range.Text = \"LTI - Java Solution\"\n\n# Paragraph 2: question 1 - change text and add a new \"b)w\" answer line\n$d.Paragraphs(2).Range.Text = \"1)adf\" + $vt + \"a)d\" + $vt + \"b)w\" + $vt + $vt\n\n# Paragraph 3: question 2 - change text and remove the \"b)3\" answer line\n$d.Paragraphs(3).Range.Text = \"2)ffff\" + $vt + \"a)1\" + $vt + $vt\n\n# Paragraph 4: question 3 - change text and remove the \"b)o\" answer line\n$d.Paragraphs(4).Range.Text = \"3)dfa\" + $vt + \"a)d\" + $vt + $vt\n"}
